# Adds the "Location" suite of sheets (Count / Names / ResultSelect) to the
# TripAdvisor Owners Page expected-elements workbook, mirroring the
# "12 Cases - TripAdvisor Owners Page" commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Location SearchResult Count" - single numeric cell
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCount = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsCount.Name = "Location SearchResult Count"
$wsCount.Range("A1").Value = 8

# ---------------------------------------------------------------------------
# 2. "Location SearchResult Names" - 8 wrapped, multi-line location names
# ---------------------------------------------------------------------------
$wsNames = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsCount)
$wsNames.Name = "Location SearchResult Names"

$nl = [char]10
$locationNames = @(
    "New York City${nl}New York, United States",
    "New York${nl}United States, North America",
    "New York Mills${nl}Minnesota, United States",
    "New York Mills${nl}New York, United States",
    "Thousand Islands New York${nl}New York, United States",
    "West New York${nl}New Jersey, United States",
    "New Paltz${nl}New York, United States",
    "New Rochelle${nl}New York, United States"
)

$wsNames.Columns.Item(1).ColumnWidth = 28.6

for ($i = 0; $i -lt $locationNames.Length; $i++) {
    $row = $i + 1
    $cell = $wsNames.Cells.Item($row, 1)
    $cell.Value = $locationNames[$i]
    $cell.WrapText = $true
    $wsNames.Rows.Item($row).RowHeight = 30
}

# ---------------------------------------------------------------------------
# 3. "Location ResultSelect" - single selected location, becomes active sheet
# ---------------------------------------------------------------------------
$wsSelect = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsNames)
$wsSelect.Name = "Location ResultSelect"
$wsSelect.Columns.Item(1).ColumnWidth = 23.75
$wsSelect.Range("A1").Value = "New York City, New York"
$wsSelect.Range("B5").Select() | Out-Null

Write-Host "Added Location sheets; active sheet is now $($wb.ActiveSheet.Name)"
